$d = $word.ActiveDocument

# The document has one section whose header/footer slots hold three
# inline "logo" pictures:
#   - Header (first page, index 2)  -> BTec_Logo-Orange   (was name "image2.jpg", -> "image1.jpg")
#   - Footer (default, index 1)     -> PearsonLogo.png    (was name "image1.png", -> "image2.png")
#   - Footer (first page, index 2)  -> PearsonLogo.png    (was name "image1.png", -> "image2.png")
#
# This renames each inline picture's display "Name" (wp:docPr/@name &
# mirrored pic:cNvPr/@name in the OOXML), swapping image1.* <-> image2.*
# without touching the image content/relationship itself.

$sec = $d.Sections.Item(1)

# --- Headers: rename the BTec logo from image2.jpg -> image1.jpg ---
for ($h = 1; $h -le 3; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}

# --- Footers: rename both Pearson logos from image1.png -> image2.png ---
for ($f = 1; $f -le 3; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}
